$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Title in A1: "13 metabolites" -> "11 metabolites", keep the rich-text runs ---
$titleCell = $ws.Range("A1")
$titleCell.Characters(50, 2).Text = "11"

$titleText = $titleCell.Value2
$italicWord = "Bordetella pertussis "
$idx1 = $titleText.IndexOf($italicWord)
$run1 = $titleCell.Characters($idx1 + 1, $italicWord.Length)
$run1.Font.Italic = $true

$tailWord = "and were not in iBP1870"
$idx2 = $titleText.IndexOf($tailWord)
$run2 = $titleCell.Characters($idx2 + 1, $tailWord.Length)
$run2.Font.Italic = $false

# --- 2. Update the metabolic-process rows (new total = 22) ---
# Row 4: Carbohydrate metabolism
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 36.363636363636367
$ws.Range("C4").NumberFormat = "0.0"

# Row 5: Glycan biosynthesis and metabolism
$ws.Range("C5").Value = 18.181818181818183
$ws.Range("C5").NumberFormat = "0.0"

# Row 6: Biosynthesis of other secondary metabolites
$ws.Range("C6").Value = 18.181818181818183
$ws.Range("C6").NumberFormat = "0.0"

# Row 7: Metabolism of terpenoids and polyketides
$ws.Range("C7").Value = 13.636363636363637
$ws.Range("C7").NumberFormat = "0.0"

# Row 8: now "ABC transporters" (was "Xenobiotics biodegradation and metabolism")
$ws.Range("A8").Value = "ABC transporters"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 9.0909090909090917
$ws.Range("C8").NumberFormat = "0.0"

# Row 9: now "Phosphotransferase system (PTS)" (was "ABC transporters")
$ws.Range("A9").Value = "Phosphotransferase system (PTS)"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 4.5454545454545459
$ws.Range("C9").NumberFormat = "0.0"

# --- 3. Remove the trailing rows (Nucleotide metabolism / Amino acid metabolism /
#        Metabolism of other amino acids), keeping their formatting intact ---
$ws.Range("A10:B13").ClearContents()
$ws.Range("C10:C13").ClearContents()
